# Custom avlo to increase transp motor gas use
#
# Primary change: the manual "Average Vehicle Loading" input on the
# AVLo-passengers sheet (cell B2) is updated from 1.67 to 1.5. All of
# C2:AK2 are a shared formula (=$B2) that recalculates automatically.
#
# A couple of incidental row-height resets (rows whose cached/auto wrap
# height no longer needs an explicit override) are also reproduced.

$wb = $excel.ActiveWorkbook

# --- AVLo-passengers: update the core input value -----------------------
$wsPax = $wb.Worksheets.Item("AVLo-passengers")
$wsPax.Activate()
$wsPax.Range("B2").Select()
$wsPax.Range("B2").Value = 1.5

# Row 1 (header row, wrapped text) settles to a shorter auto height.
$wsPax.Rows.Item(1).RowHeight = 45

# --- AVLo-freight: header row height also settles shorter ---------------
$wsFreight = $wb.Worksheets.Item("AVLo-freight")
$wsFreight.Rows.Item(1).RowHeight = 45

# --- BTS NTS Modal Profile Data: a few rows drop their custom height ----
$wsBts = $wb.Worksheets.Item("BTS NTS Modal Profile Data")
$wsBts.Rows.Item(36).AutoFit()
$wsBts.Rows.Item(37).AutoFit()
$wsBts.Rows.Item(60).AutoFit()
